$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5,1).Value = "Globo"
$ws.Cells.Item(5,2).Value = "Bom Dia Inter"
$ws.Cells.Item(5,3).Value = "Limpeza Pública"
$ws.Cells.Item(5,4).Value = "2025-04-02T16:22"
$ws.Cells.Item(5,5).Value = "Neutro"
$ws.Cells.Item(5,6).Value = "aaaa teste"

$ws.Cells.Item(6,1).Value = "Globo"
$ws.Cells.Item(6,2).Value = "Bom Dia Inter"
$ws.Cells.Item(6,3).Value = "Social"
$ws.Cells.Item(6,4).Value = "2025-04-02T08:04"
$ws.Cells.Item(6,5).Value = "Neutro"
$ws.Cells.Item(6,6).Value = "Mutirão CadÚnico em Campos. Serão distribuídas 800 senhas. Recadastramento accontece no Parque Turf Club. Repórter *ao vivo* no local. Muitas pessoas em frente à unidade. Filas enormes antes de abrir `nos portões. Conteúdo informativo, mas muitas reclamações de assistida, incluindo do atendimento e do tempo de espera. "

$ws.Cells.Item(7,1).Value = "Globo"
$ws.Cells.Item(7,2).Value = "RJ TV 1"
$ws.Cells.Item(7,3).Value = "PROCON"
$ws.Cells.Item(7,4).Value = "2025-04-02T12:51"
$ws.Cells.Item(7,5).Value = "Positivo"
$ws.Cells.Item(7,6).Value = "Variação de preços dos produtos da Páscoa. Em Canpos, a fiscalização encontrou 75% de variação no preço de alguns produtos. Repórter *ao vivo*"

Write-Output "done"
